$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row: rename the "<name>_old" / "<name>_new" column headers to use
#    the respective input-file-version suffixes "_FV2404" / "_FV2410".
# ---------------------------------------------------------------------------
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    if ($v -like "*_old") {
        $cell.Value = ($v -replace "_old$", "_FV2404")
    } elseif ($v -like "*_new") {
        $cell.Value = ($v -replace "_new$", "_FV2410")
    }
}

# ---------------------------------------------------------------------------
# 2) Turn the populated range into a real Excel Table ("Table1") so the
#    header row gets filter buttons / structured references, matching the
#    workbook's new xlsx export format.
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:U79")
$tbl = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# ---------------------------------------------------------------------------
# 3) Freeze the header row so it stays visible while scrolling.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
